# ArchitectureDiagram.pptx edit: update the "Message system" portion of the
# diagram (remove the old Message/MessageHandler shapes + their connectors
# and labels from slides 1 and 2) and refresh the slide master / layout
# "last edited" date placeholders that PowerPoint stamps on every save.

function Remove-ShapeById($shapes, $id) {
    for ($i = $shapes.Count; $i -ge 1; $i--) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            $sh.Delete()
            return
        }
    }
}

$p = $ppt.ActivePresentation

# --- 1. Refresh the datetimeFigureOut date placeholder everywhere it lives:
#        the slide master and all of its custom layouts. ---
$newDate = "26/10/2020"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "27/5/2020") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "27/5/2020") {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Slide 1: remove the "MessageHandler" rounded rectangle, its dashed
#        arrow from the Output rectangle, the "Polls for Messages" label,
#        the "Message" rounded rectangle near the top, the elbow connector
#        and dashed arrow feeding it, and the "Sends Messages" label. ---
$slide1 = $p.Slides.Item(1)
Remove-ShapeById $slide1.Shapes 30   # Rectangle: Rounded Corners 29 ("MessageHandler")
Remove-ShapeById $slide1.Shapes 31   # Straight Arrow Connector 30
Remove-ShapeById $slide1.Shapes 33   # TextBox 32 ("Polls for Messages")
Remove-ShapeById $slide1.Shapes 60   # Rectangle: Rounded Corners 59 ("Message")
Remove-ShapeById $slide1.Shapes 83   # Connector: Elbow 82
Remove-ShapeById $slide1.Shapes 88   # Straight Arrow Connector 87
Remove-ShapeById $slide1.Shapes 91   # TextBox 90 ("Sends Messages")

# --- 3. Slide 2: remove the leftover "Messages" arrow/label wiring that
#        pointed at the now-deleted Message shapes. ---
$slide2 = $p.Slides.Item(2)
Remove-ShapeById $slide2.Shapes 71    # Straight Arrow Connector 70
Remove-ShapeById $slide2.Shapes 85    # TextBox 84 ("Messages")
Remove-ShapeById $slide2.Shapes 104   # Straight Arrow Connector 103
Remove-ShapeById $slide2.Shapes 106   # Straight Arrow Connector 105
